$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C. This pushes the existing
# column C (per-report rank/history data) out to column E, leaving
# fresh blank columns C and D behind - exactly mirroring how a new
# "as of" rank column gets inserted every refresh.
$ws.Columns("C:D").Insert()

# Shift the header dates: the old B1 ("Jun_13") and C1 ("Jun_10")
# values move right into D1/E1 (already preserved by the column
# insert above for C1->E1; B1 still needs to move to D1), and the two
# new leading columns get the newest dates.
$ws.Range("D1").Value = "Jun_13"
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# The two freshly inserted columns don't have rank data yet for any
# company, so fill them with the same "unranked" placeholder used
# elsewhere in column B.
$ws.Range("C2:D27").Value = "UN"

# Match the original column C formatting (width 8, explicit custom
# width) on the new C/D columns and keep it on the shifted-out E.
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.14
